# The "Xcg ESTIMATION METHOD COMPARISON (WITHOUT CALIBRATIONS)" tables on the
# FUSELAGE and WING sheets list the TORENBEEK_1982 and SFORZA estimates.
# Swap the row order of these two methods (SFORZA now comes first), keeping
# each method's own value attached to its label.

$wb = $excel.ActiveWorkbook

$fuselage = $wb.Worksheets.Item("FUSELAGE")
$fuselage.Range("A23").Value = "SFORZA"
$fuselage.Range("C23").Value = 17.143322222222217
$fuselage.Range("A24").Value = "TORENBEEK_1982"
$fuselage.Range("C24").Value = 16.8345

$wing = $wb.Worksheets.Item("WING")
$wing.Range("A23").Value = "SFORZA"
$wing.Range("C23").Value = 4.3631082000119275
$wing.Range("A24").Value = "TORENBEEK_1982"
$wing.Range("C24").Value = 3.5939754358446514

$wing.Range("A27").Value = "SFORZA"
$wing.Range("C27").Value = 4.998846772296348
$wing.Range("A28").Value = "TORENBEEK_1982"
$wing.Range("C28").Value = 6.114221148470394
